$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, shifting rows 172:259 down to 173:260
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new record's data
$ws.Cells.Item(172, 1).Value = 8
$ws.Cells.Item(172, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(172, 3).Value = "Coquimbo"
$ws.Cells.Item(172, 4).Value = 44582
$ws.Cells.Item(172, 5).Value = 4
$ws.Cells.Item(172, 6).Value = 100112043
$ws.Cells.Item(172, 7).Value = "Pepino dulce"
$ws.Cells.Item(172, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 500
$ws.Cells.Item(172, 11).Value = 16000
$ws.Cells.Item(172, 12).Value = 17000
$ws.Cells.Item(172, 13).Value = 16500
$ws.Cells.Item(172, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(172, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(172, 16).Value = 917
$ws.Cells.Item(172, 17).Value = 18
$ws.Cells.Item(172, 18).Value = "Hortaliza"
